# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt -
# Cebollín" right after the existing row 126, pushing the following rows
# (old 127..192) down to 128..193, and fill the newly opened row 127 with
# the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 127..192 down to 128..193, leaving row 127 empty.
$ws.Rows.Item(127).Insert()

# Populate the new row 127 with the new weekly record.
$ws.Range("A127").Value = 4
$ws.Range("B127").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C127").Value = "Los Lagos"
$ws.Range("D127").Value = 44523
$ws.Range("E127").Value = 10
$ws.Range("F127").Value = 100112037
$ws.Range("G127").Value = "Cebollín"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 180
$ws.Range("K127").Value = 6000
$ws.Range("L127").Value = 6000
$ws.Range("M127").Value = 6000
$ws.Range("N127").Value = "$/paquete 36 unidades"
$ws.Range("O127").Value = "Región Metropolitana"
$ws.Range("P127").Value = 167
$ws.Range("Q127").Value = 36
$ws.Range("R127").Value = "Hortaliza"
